# Trade #129 closed at 2026-02-16 21:48:59 - leadlag UP +0.000%
# This script applies the bookkeeping updates that follow from trade #129
# being opened (leadlag, UP) and several earlier trades (#109-#116) being
# closed out (time_exit_5min), updating every sheet that tracks them.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: write a value to a cell while forcing literal text for anything
# that Excel would otherwise "smart convert" (dates, times, %, numbers-as-
# text). A leading apostrophe marks the entry as text, matching how Excel
# itself preserves user-typed text, and it keeps the default "General"
# number format / style intact (no extra style records are introduced).
# ---------------------------------------------------------------------------

# =================== Sheet: leadlag ===================
$ws = $wb.Worksheets.Item("leadlag")

# Rows 85-91 (trades 109-111, 113-116) change from OPEN -> CLOSED
$ws.Range("G85").Value = 67511.163512
$ws.Range("H85").Value = "CLOSED"
$ws.Range("I85").Value = 1.2932
$ws.Range("J85").Value = 12.93
$ws.Range("M85").Value = "time_exit_5min"
$ws.Range("N85").Value = 5

$ws.Range("G86").Value = 67908.535476
$ws.Range("H86").Value = "CLOSED"
$ws.Range("I86").Value = 0.6209
$ws.Range("J86").Value = 6.21
$ws.Range("M86").Value = "time_exit_5min"
$ws.Range("N86").Value = 5

$ws.Range("G87").Value = 68082.997787
$ws.Range("H87").Value = "CLOSED"
$ws.Range("I87").Value = 0.3943
$ws.Range("J87").Value = 3.94
$ws.Range("M87").Value = "time_exit_5min"
$ws.Range("N87").Value = 5

$ws.Range("G88").Value = 68476.49535300001
$ws.Range("H88").Value = "CLOSED"
$ws.Range("I88").Value = 0.1055
$ws.Range("J88").Value = 1.06
$ws.Range("M88").Value = "time_exit_5min"
$ws.Range("N88").Value = 5

$ws.Range("G89").Value = 68905.925091
$ws.Range("H89").Value = "CLOSED"
$ws.Range("I89").Value = 0.6977
$ws.Range("J89").Value = 6.98
$ws.Range("M89").Value = "time_exit_5min"
$ws.Range("N89").Value = 5

$ws.Range("G90").Value = 68618.97139399999
$ws.Range("H90").Value = "CLOSED"
$ws.Range("I90").Value = 0.2652
$ws.Range("J90").Value = 2.65
$ws.Range("M90").Value = "time_exit_5min"
$ws.Range("N90").Value = 5

$ws.Range("G91").Value = 67876.50332
$ws.Range("H91").Value = "CLOSED"
$ws.Range("I91").Value = 0.7494
$ws.Range("J91").Value = 7.49
$ws.Range("M91").Value = "time_exit_5min"
$ws.Range("N91").Value = 5

# New row 104: trade #129 just opened (leadlag, UP)
$ws.Range("A104").Value = 129
$ws.Range("B104").Value = "'2026-02-16"
$ws.Range("C104").Value = "'21:48:59"
$ws.Range("D104").Value = "leadlag"
$ws.Range("E104").Value = "UP"
$ws.Range("F104").Value = 68376.295
$ws.Range("H104").Value = "OPEN"
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("K104").Value = 0.75
$ws.Range("L104").Value = "Coinbase leading with 0.095% move"
$ws.Range("N104").Value = 0

# =================== Sheet: momentum ===================
$ws = $wb.Worksheets.Item("momentum")

# Row 27 (trade 112) changes from OPEN -> CLOSED
$ws.Range("G27").Value = 68011.664288
$ws.Range("H27").Value = "CLOSED"
$ws.Range("I27").Value = 0.5034
$ws.Range("J27").Value = 5.03
$ws.Range("M27").Value = "time_exit_5min"
$ws.Range("N27").Value = 5

# =================== Sheet: All Trades ===================
$ws = $wb.Worksheets.Item("All Trades")

# New rows 110-117, combining the newly-closed trades 109-116
# (leadlag + momentum), in trade-number order.

# Row 110 - trade 109 (leadlag)
$ws.Range("A110").Value = 109
$ws.Range("B110").Value = "'2026-02-16"
$ws.Range("C110").Value = "'21:42:45"
$ws.Range("D110").Value = "leadlag"
$ws.Range("E110").Value = "DOWN"
$ws.Range("F110").Value = 68395.645
$ws.Range("G110").Value = 67511.163512
$ws.Range("H110").Value = "CLOSED"
$ws.Range("I110").Value = 1.2932
$ws.Range("J110").Value = 12.93
$ws.Range("K110").Value = 0.75
$ws.Range("L110").Value = "Binance leading with -0.107% move"
$ws.Range("M110").Value = "time_exit_5min"
$ws.Range("N110").Value = 5

# Row 111 - trade 110 (leadlag)
$ws.Range("A111").Value = 110
$ws.Range("B111").Value = "'2026-02-16"
$ws.Range("C111").Value = "'21:42:52"
$ws.Range("D111").Value = "leadlag"
$ws.Range("E111").Value = "DOWN"
$ws.Range("F111").Value = 68332.81
$ws.Range("G111").Value = 67908.535476
$ws.Range("H111").Value = "CLOSED"
$ws.Range("I111").Value = 0.6209
$ws.Range("J111").Value = 6.21
$ws.Range("K111").Value = 0.75
$ws.Range("L111").Value = "Coinbase leading with -0.113% move"
$ws.Range("M111").Value = "time_exit_5min"
$ws.Range("N111").Value = 5

# Row 112 - trade 111 (leadlag)
$ws.Range("A112").Value = 111
$ws.Range("B112").Value = "'2026-02-16"
$ws.Range("C112").Value = "'21:42:58"
$ws.Range("D112").Value = "leadlag"
$ws.Range("E112").Value = "DOWN"
$ws.Range("F112").Value = 68352.48
$ws.Range("G112").Value = 68082.997787
$ws.Range("H112").Value = "CLOSED"
$ws.Range("I112").Value = 0.3943
$ws.Range("J112").Value = 3.94
$ws.Range("K112").Value = 0.6602
$ws.Range("L112").Value = "Coinbase leading with -0.066% move"
$ws.Range("M112").Value = "time_exit_5min"
$ws.Range("N112").Value = 5

# Row 113 - trade 112 (momentum)
$ws.Range("A113").Value = 112
$ws.Range("B113").Value = "'2026-02-16"
$ws.Range("C113").Value = "'21:43:05"
$ws.Range("D113").Value = "momentum"
$ws.Range("E113").Value = "DOWN"
$ws.Range("F113").Value = 68355.74000000001
$ws.Range("G113").Value = 68011.664288
$ws.Range("H113").Value = "CLOSED"
$ws.Range("I113").Value = 0.5034
$ws.Range("J113").Value = 5.03
$ws.Range("K113").Value = 0.9
$ws.Range("L113").Value = "Downward momentum: -0.175% over 10 samples"
$ws.Range("M113").Value = "time_exit_5min"
$ws.Range("N113").Value = 5

# Row 114 - trade 113 (leadlag)
$ws.Range("A114").Value = 113
$ws.Range("B114").Value = "'2026-02-16"
$ws.Range("C114").Value = "'21:43:17"
$ws.Range("D114").Value = "leadlag"
$ws.Range("E114").Value = "UP"
$ws.Range("F114").Value = 68404.295
$ws.Range("G114").Value = 68476.49535300001
$ws.Range("H114").Value = "CLOSED"
$ws.Range("I114").Value = 0.1055
$ws.Range("J114").Value = 1.06
$ws.Range("K114").Value = 0.75
$ws.Range("L114").Value = "Coinbase leading with 0.078% move"
$ws.Range("M114").Value = "time_exit_5min"
$ws.Range("N114").Value = 5

# Row 115 - trade 114 (leadlag)
$ws.Range("A115").Value = 114
$ws.Range("B115").Value = "'2026-02-16"
$ws.Range("C115").Value = "'21:43:23"
$ws.Range("D115").Value = "leadlag"
$ws.Range("E115").Value = "UP"
$ws.Range("F115").Value = 68428.50999999999
$ws.Range("G115").Value = 68905.925091
$ws.Range("H115").Value = "CLOSED"
$ws.Range("I115").Value = 0.6977
$ws.Range("J115").Value = 6.98
$ws.Range("K115").Value = 0.75
$ws.Range("L115").Value = "Binance leading with 0.077% move"
$ws.Range("M115").Value = "time_exit_5min"
$ws.Range("N115").Value = 5

# Row 116 - trade 115 (leadlag)
$ws.Range("A116").Value = 115
$ws.Range("B116").Value = "'2026-02-16"
$ws.Range("C116").Value = "'21:43:30"
$ws.Range("D116").Value = "leadlag"
$ws.Range("E116").Value = "UP"
$ws.Range("F116").Value = 68437.495
$ws.Range("G116").Value = 68618.97139399999
$ws.Range("H116").Value = "CLOSED"
$ws.Range("I116").Value = 0.2652
$ws.Range("J116").Value = 2.65
$ws.Range("K116").Value = 0.6051
$ws.Range("L116").Value = "Binance leading with 0.061% move"
$ws.Range("M116").Value = "time_exit_5min"
$ws.Range("N116").Value = 5

# Row 117 - trade 116 (leadlag)
$ws.Range("A117").Value = 116
$ws.Range("B117").Value = "'2026-02-16"
$ws.Range("C117").Value = "'21:43:42"
$ws.Range("D117").Value = "leadlag"
$ws.Range("E117").Value = "DOWN"
$ws.Range("F117").Value = 68388.985
$ws.Range("G117").Value = 67876.50332
$ws.Range("H117").Value = "CLOSED"
$ws.Range("I117").Value = 0.7494
$ws.Range("J117").Value = 7.49
$ws.Range("K117").Value = 0.75
$ws.Range("L117").Value = "Binance leading with -0.079% move"
$ws.Range("M117").Value = "time_exit_5min"
$ws.Range("N117").Value = 5

# =================== Sheet: Summary ===================
$ws = $wb.Worksheets.Item("Summary")

$ws.Range("C2").Value = 116
$ws.Range("D2").Value = "'74.1%"
$ws.Range("E2").Value = "'+39.5445%"
$ws.Range("F2").Value = "'+0.3409%"

$ws.Range("C3").Value = 102
$ws.Range("D3").Value = "'60.8%"
$ws.Range("E3").Value = "'+25.2305%"
$ws.Range("F3").Value = "'+0.2474%"

$ws.Range("D4").Value = "'92.3%"
$ws.Range("E4").Value = "'+14.3140%"
$ws.Range("F4").Value = "'+0.5505%"

# =================== Sheet: Comparison ===================
$ws = $wb.Worksheets.Item("Comparison")

$ws.Range("B2").Value = 102
$ws.Range("C2").Value = "'60.8%"
$ws.Range("D2").Value = "'4.13"
$ws.Range("E2").Value = "'+0.5370%"
$ws.Range("G2").Value = "'1.86"

$ws.Range("C3").Value = "'92.3%"
$ws.Range("D3").Value = "'13.73"
$ws.Range("E3").Value = "'+0.6433%"
$ws.Range("G3").Value = "'1.14"
